$d = $word.ActiveDocument

# 1. Replace the two runs' text ("**ID__AFFARS_pgi_5301_topic_34__ID**" + " ")
#    with a single run containing the new bookmark-style placeholder text.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5301_topic_34__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5301_9001_92__ID**", 2)

# 2. Update the first paragraph's formatting: add a paragraph border (5pt/twip
#    spacing on all sides) and change the left indent from 120 to 225 twips
#    (i.e. from 6pt to 11.25pt).
$p1 = $d.Paragraphs.Item(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25
